$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the two pairs of adjacent country rows (labels only; stats stay with their row) ---
# Fiyi / Santa Lucia / Timor Oriental  ->  Fiyi / Timor Oriental / Santa Lucia
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# Groenlandia / Montserrat / Islas Malvinas  ->  Groenlandia / Islas Malvinas / Montserrat
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 23:49"

# --- Update numeric data cells ---
# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 6997288
$ws.Range("C4").Value = 29885
$ws.Range("D4").Value = 4244789
$ws.Range("E4").Value = 2548416
$ws.Range("G4").Value = 259
$ws.Range("H4").Value = 204083

# Row 6 (Rusia)
$ws.Range("B6").Value = 4544629
$ws.Range("C6").Value = 16282
$ws.Range("E6").Value = 587639
$ws.Range("G6").Value = 330
$ws.Range("H6").Value = 136895

# Row 54
$ws.Range("B54").Value = 65039
$ws.Range("C54").Value = 540
$ws.Range("D54").Value = 57950
$ws.Range("E54").Value = 6868

# Row 84
$ws.Range("B84").Value = 18863
$ws.Range("C84").Value = 44
$ws.Range("D84").Value = 13580
$ws.Range("E84").Value = 4522
$ws.Range("G84").Value = 6
$ws.Range("H84").Value = 761

# Row 157
$ws.Range("B157").Value = 1666
$ws.Range("C157").Value = 7
$ws.Range("D157").Value = 1269
$ws.Range("E157").Value = 356

# Row 166
$ws.Range("B166").Value = 1151
$ws.Range("C166").Value = 2
$ws.Range("D166").Value = 967
$ws.Range("E166").Value = 103

# Row 214 (now "Islas Malvinas" after the label swap above)
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215 (now "Montserrat" after the label swap above)
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$wb.Save()
